$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# Fix title: "МатематическиЯ" -> "МатематическаЯ"
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "МатематическаЯ и физическая модели"

# Fix content placeholder first run: "Основыне" -> "Основные"
$contentShape = $s.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$firstRun = $tr.Characters(1, 29)
$firstRun.Text = "Основные формулы при создании"
